$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ------------------------------------------------------------------
# 1) Update the date in A1 (2024-04-24 -> 2024-05-24, serial 45406 -> 45436)
# ------------------------------------------------------------------
$ws.Range("A1").Value = 45436

# ------------------------------------------------------------------
# 2) Update the price in D44
# ------------------------------------------------------------------
$ws.Range("D44").Value = 43783.243

# ------------------------------------------------------------------
# 3) Re-sequence the merged cell ranges so the <mergeCells> list comes
#    back out in the target order. The engine appends a range to the
#    end of the merge list whenever it is (re)merged, so re-merging in
#    the desired final order reproduces that order exactly.
#
#    Re-merging the already-populated label ranges (B42:C42, B43:C43,
#    B44:C44) also nudges Excel's automatic border adjustment on the
#    cells inside them, so their formatting is stashed in a few unused
#    cells first and restored afterwards to keep the sheet visually
#    identical to before.
# ------------------------------------------------------------------

$stashPairs = @(
    ,@("B42", "B41")
    ,@("C42", "C41")
    ,@("B43", "D41")
    ,@("C43", "B46")
    ,@("B44", "C46")
    ,@("C44", "D46")
)

foreach ($pair in $stashPairs) {
    $ws.Range($pair[0]).Copy() | Out-Null
    $ws.Range($pair[1]).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = $false

$order = @("A1:D1", "B44:C44", "A9:D9", "B42:C42", "A11:D11", "A10:D10", "B43:C43")
foreach ($r in $order) {
    $ws.Range($r).UnMerge()
    $ws.Range($r).Merge()
}

foreach ($pair in $stashPairs) {
    $ws.Range($pair[1]).Copy() | Out-Null
    $ws.Range($pair[0]).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range($pair[1]).Clear() | Out-Null
}
$excel.CutCopyMode = $false
